$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows 72 & 73: start_time (D) changes from 9 -> 13,
# which also affects the computed end_time (F, via SUM(D,E)) and the
# manually-assigned time_slot (H) which moves from 16 -> 18 to avoid two
# consecutive timeslots for the same invigilator in a day.
$ws.Cells.Item(72, 4).Value = 13
$ws.Cells.Item(72, 8).Value = 18

$ws.Cells.Item(73, 4).Value = 13
$ws.Cells.Item(73, 8).Value = 18

# --- Add the new exam row (row 76) that was missing previously, which is
# the fix for "the last exam had no colour/time_slot assigned to it".
$ws.Cells.Item(76, 1).Formula = "=ROW()-2"

$ws.Cells.Item(76, 2).Value = "Exam75"

$ws.Cells.Item(76, 3).Value2 = 45437
$ws.Cells.Item(76, 3).NumberFormat = $ws.Cells.Item(75, 3).NumberFormat

$ws.Cells.Item(76, 4).Value = 13
$ws.Cells.Item(76, 5).Value = 2

$ws.Cells.Item(76, 6).Formula = "=SUM(D76,E76)"

$ws.Cells.Item(76, 7).Value = 46

$ws.Cells.Item(76, 8).Value2 = 18
$ws.Cells.Item(76, 8).NumberFormat = $ws.Cells.Item(75, 8).NumberFormat

# Recalculate so the formula-driven cells (A76, F76, and the existing
# shared formulas in F72/F73) reflect the updated values.
$wb.Application.Calculate()

# Bring the new last row into view, matching how Excel scrolls/selects
# after the row is appended.
$ws.Activate() | Out-Null
$ws.Range("G76").Select() | Out-Null
